$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append three new monthly rows (275-277) below the existing series ---
$ws.Range("A275").Value = 44835
$ws.Range("A276").Value = 44866
$ws.Range("A277").Value = 44896

$ws.Range("B275").Formula = "=B274+_xlfn.STDEV.S(B272:B274)"
$ws.Range("B276").Formula = "=B275+_xlfn.STDEV.S(B273:B275)"
$ws.Range("B277").Formula = "=B276+_xlfn.STDEV.S(B274:B276)"

# --- Highlight the new study rows (and the two rows feeding it) with a
#     yellow fill, keeping the existing date number format on column A ---
$ws.Range("A273:A277").Interior.Color = 65535
$ws.Range("A273:A277").NumberFormat = "[$-409]mmm\-yy;@"
$ws.Range("B273:B277").Interior.Color = 65535

# --- Move the active selection to D274, matching the author's last click ---
$ws.Range("D274").Select() | Out-Null
